$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add Wins, Losses, Ties in columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match header style (bold, centered top, thin border) used by the rest of row 1
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Data rows 2-48: Wins = 65, Losses = 97, Ties = 0 for every row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 29).Value = 65
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 0
}
